# Update vm_pu results for Case 1_70 (380 kV case)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    "B2" = 1.02; "C2" = 1.027772848844516; "D2" = 1.03241419393604; "E2" = 1.036476268768222; "F2" = 1.044259761674375; "I2" = 1.034809181180663; "J2" = 1.032928861695017; "K2" = 1.035219493614367; "L2" = 1.039269901277221; "M2" = 1.047031313507817; "N2" = 1.014996663321528
    "B3" = 1.02; "C3" = 1.028616961481953; "D3" = 1.032872608063823; "E3" = 1.037260420819911; "F3" = 1.045205447183924; "I3" = 1.034955026496421; "J3" = 1.033413842314216; "K3" = 1.035487612770188; "L3" = 1.039863730229342; "M3" = 1.047787849765716; "N3" = 1.015157899109492
    "B4" = 1.02; "C4" = 1.029163582686708; "D4" = 1.033169167823498; "E4" = 1.037768653031377; "F4" = 1.045818547608878; "I4" = 1.035048005445274; "J4" = 1.033727434590059; "K4" = 1.035660312548436; "L4" = 1.040248153200057; "M4" = 1.048277933184906; "N4" = 1.015262127956719
    "B5" = 1.02; "C5" = 1.029393482237638; "D5" = 1.033293824113352; "E5" = 1.037982511688944; "F5" = 1.046076575488493; "I5" = 1.035086759504164; "J5" = 1.033859214299629; "K5" = 1.035732724825057; "L5" = 1.040409805400067; "M5" = 1.048484095337457; "N5" = 1.015305921036701
    "B6" = 1.02; "C6" = 1.029432089204633; "D6" = 1.033314753374892; "E6" = 1.038018431050853; "F6" = 1.046119915901433; "I6" = 1.035093246852837; "J6" = 1.033881337463225; "K6" = 1.035744871944146; "L6" = 1.040436949868474; "M6" = 1.048518718512773; "N6" = 1.01531327262497
    "B7" = 1.02; "C7" = 1.029166654224695; "D7" = 1.033170833557864; "E7" = 1.03777150984598; "F7" = 1.045821994290162; "I7" = 1.035048524593534; "J7" = 1.033729195652562; "K7" = 1.035661280875803; "L7" = 1.040250313046092; "M7" = 1.048280687420512; "N7" = 1.0152627132197
    "B8" = 1.02; "C8" = 1.028058032002591; "D8" = 1.032569129414199; "E8" = 1.036741103053617; "F8" = 1.044579115438718; "I8" = 1.034858758137451; "J8" = 1.03309280854505; "K8" = 1.035310268799668; "L8" = 1.039470551321519; "M8" = 1.047286872883848; "N8" = 1.015051174466124
    "B9" = 1.02; "C9" = 1.026107809944087; "D9" = 1.031508437500422; "E9" = 1.034931851745287; "F9" = 1.04239810700035; "I9" = 1.034513737184858; "J9" = 1.031969760662109; "K9" = 1.034685737287288; "L9" = 1.038097917901567; "M9" = 1.045539953623766; "N9" = 1.014677657871808
    "B10" = 1.02; "C10" = 1.024809976376596; "D10" = 1.030801145996034; "E10" = 1.033730125511061; "F10" = 1.040950323381812; "I10" = 1.034276629016784; "J10" = 1.031220021493878; "K10" = 1.0342654228818; "L10" = 1.037183853731219; "M10" = 1.044378327774678; "N10" = 1.014428162516623
    "B11" = 1.02; "C11" = 1.02424856737909; "D11" = 1.03049486620263; "E11" = 1.033210839209928; "F11" = 1.040324914207925; "I11" = 1.0341722867394; "J11" = 1.03089514411232; "K11" = 1.034082498906558; "L11" = 1.036788313539685; "M11" = 1.043876058288975; "N11" = 1.014320018762774
    "B12" = 1.02; "C12" = 1.024040120890731; "D12" = 1.030381099508879; "E12" = 1.033018115466568; "F12" = 1.040092835066327; "I12" = 1.03413327889193; "J12" = 1.030774436000151; "K12" = 1.034014415105991; "L12" = 1.036641432035025; "M12" = 1.043689603045707; "N12" = 1.014279833165258
    "B13" = 1.02; "C13" = 1.024084829463673; "D13" = 1.030405502853239; "E13" = 1.033059447991714; "F13" = 1.040142606590366; "I13" = 1.034141657528723; "J13" = 1.030800329828167; "K13" = 1.034029025522518; "L13" = 1.036672936797002; "M13" = 1.043729593326032; "N13" = 1.014288453838222
    "B14" = 1.02; "C14" = 1.024231335368438; "D14" = 1.030485462216124; "E14" = 1.033194905284325; "F14" = 1.040305725860131; "I14" = 1.034169067443355; "J14" = 1.030885167032736; "K14" = 1.034076873879793; "L14" = 1.036776171450446; "M14" = 1.043860643587556; "N14" = 1.014316697336596
    "B15" = 1.02; "C15" = 1.02432161390577; "D15" = 1.030534727778472; "E15" = 1.033278386551201; "F15" = 1.040406258980832; "I15" = 1.034185922430823; "J15" = 1.030937433539406; "K15" = 1.034106336624914; "L15" = 1.036839783032478; "M15" = 1.043941402601064; "N15" = 1.014334096955231
    "B16" = 1.02; "C16" = 1.024847247196788; "D16" = 1.03082147256653; "E16" = 1.033764611501533; "F16" = 1.04099186128221; "I16" = 1.034283518697185; "J16" = 1.03124157767697; "K16" = 1.034277543563697; "L16" = 1.037210109930252; "M16" = 1.044411677090403; "N16" = 1.014435337367254
    "B17" = 1.02; "C17" = 1.025177114426662; "D17" = 1.031001336581337; "E17" = 1.034069894986957; "F17" = 1.041359594442688; "I17" = 1.034344290981414; "J17" = 1.031432297072347; "K17" = 1.034384690398811; "L17" = 1.037442475586061; "M17" = 1.044706862264334; "N17" = 1.014498813488832
    "B18" = 1.02; "C18" = 1.025369574480304; "D18" = 1.031106246360416; "E18" = 1.034248064717185; "F18" = 1.041574230688373; "I18" = 1.034379576961378; "J18" = 1.031543517680019; "K18" = 1.034447097993535; "L18" = 1.037578035116905; "M18" = 1.044879108360109; "N18" = 1.014535827334923
    "B19" = 1.02; "C19" = 1.025435207513673; "D19" = 1.031142017509452; "E19" = 1.034308833416449; "F19" = 1.041647440425863; "I19" = 1.034391581155457; "J19" = 1.031581437131885; "K19" = 1.034468362177435; "L19" = 1.037624261541678; "M19" = 1.044937851583256; "N19" = 1.014548446271726
    "B20" = 1.02; "C20" = 1.025141717190557; "D20" = 1.030982039052842; "E20" = 1.034037130288206; "F20" = 1.041320125277296; "I20" = 1.034337787385419; "J20" = 1.031411837034849; "K20" = 1.034373203789121; "L20" = 1.037417542398757; "M20" = 1.044675184481098; "N20" = 1.01449200420556
    "B21" = 1.02; "C21" = 1.024188190650611; "D21" = 1.030461916181553; "E21" = 1.033155011993704; "F21" = 1.040257685057415; "I21" = 1.034161002814315; "J21" = 1.030860185524814; "K21" = 1.034062787521798; "L21" = 1.036745770310543; "M21" = 1.04382204950302; "N21" = 1.014308380773835
    "B22" = 1.02; "C22" = 1.023589167859819; "D22" = 1.030134890957142; "E22" = 1.032601329123402; "F22" = 1.039590993239998; "I22" = 1.034048402559022; "J22" = 1.03051314310817; "K22" = 1.033866820488014; "L22" = 1.036323631081691; "M22" = 1.043286286896923; "N22" = 1.014192835961965
    "B23" = 1.02; "C23" = 1.023906673562217; "D23" = 1.030308252883729; "E23" = 1.032894757231604; "F23" = 1.039944294771631; "I23" = 1.034108231128614; "J23" = 1.030697135161333; "K23" = 1.033970781343597; "L23" = 1.036547392777377; "M23" = 1.043570243886726; "N23" = 1.014254097170602
    "B24" = 1.02; "C24" = 1.025157711510954; "D24" = 1.030990758781311; "E24" = 1.034051934927063; "F24" = 1.041337959250347; "I24" = 1.034340726579504; "J24" = 1.031421082115649; "K24" = 1.034378394369816; "L24" = 1.037428808557082; "M24" = 1.044689498092975; "N24" = 1.014495081060684
    "B25" = 1.02; "C25" = 1.026611587579345; "D25" = 1.03178268890264; "E25" = 1.035398811075491; "F25" = 1.042960860740722; "I25" = 1.034604187544703; "J25" = 1.032260284023784; "K25" = 1.034847897320708; "L25" = 1.038452601715854; "M25" = 1.045991054232116; "N25" = 1.014774308091626
}

foreach ($key in $data.Keys) {
    $ws.Range($key).Value = $data[$key]
}

Write-Host "Updated $($data.Count) cells"
